# Adds a new product row ("ZINCTRON 30 CAPS") to the day-sale shortage
# report, right before "املاح افونا" (which currently lives on row 41),
# and refreshes the generated-at timestamp in the footer.
#
# Strategy: insert one blank row at row 41 (this shifts every row below it
# down by one - data, merged cells, etc. all move automatically), copy the
# formatting of the row that is now directly below (row 42, which holds the
# old row-41 content/style) onto the freshly inserted row 41, then populate
# row 41 with the new product's data. Finally patch up the row that used to
# be the "total" row and the footer row (their row numbers shifted by one),
# and correct every row's height to match the final layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row right above the old row 41 ("املاح افونا") -------
$ws.Rows.Item(41).Insert()

# Bring over the cell formatting (fonts/fills/borders/number formats) from
# row 42 (the row that used to be row 41, now shifted down) so the new row
# looks exactly like every other product row.
$ws.Range("A42:Q42").Copy()
$ws.Range("A41:Q41").PasteSpecial(-4122)

# Re-create the merged cells for the new row (Insert() does not merge the
# freshly created cells automatically).
$ws.Range("A41:B41").Merge()
$ws.Range("C41:G41").Merge()
$ws.Range("H41:K41").Merge()
$ws.Range("L41:M41").Merge()
$ws.Range("N41:O41").Merge()

# --- 2. Fill in the new product's data -------------------------------------
$ws.Range("A41").Value = 35
$ws.Range("C41").Value = "ZINCTRON 30 CAPS"
$ws.Range("H41").Value = "1:1"
$ws.Range("L41").Value = "1"
$ws.Range("N41").Value = "126.00"
$ws.Range("P41").Value = "41.5800"
$ws.Range("Q41").Value = "0:1"

# --- 3. Update the running total (old row 49, now row 50) ------------------
$ws.Range("P50").Value = 1933.41

# --- 4. Update the generated-at timestamp in the footer (now row 51) -------
$ws.Range("A51").Value = "Monday, 29 September, 2025 12:56 PM"

# --- 5. Fix up row heights to match the final layout ------------------------
$ws.Rows.Item(41).RowHeight = 25.5
$ws.Rows.Item(42).RowHeight = 25.5
$ws.Rows.Item(43).RowHeight = 24.75
$ws.Rows.Item(44).RowHeight = 25.5
$ws.Rows.Item(45).RowHeight = 24.75
$ws.Rows.Item(46).RowHeight = 25.5
$ws.Rows.Item(47).RowHeight = 25.5
$ws.Rows.Item(48).RowHeight = 24.75
$ws.Rows.Item(49).RowHeight = 25.5
$ws.Rows.Item(50).RowHeight = 24.75
$ws.Rows.Item(51).RowHeight = 16.5
